$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the existing row 2 values (strings that were corrected) ---
$ws.Range("A2").Value = "Rakesh11"
$ws.Range("D2").Value = "AutomationModel"
$ws.Range("F2").Value = "AutomationAsset"
$ws.Range("H2").Value = "AutomationSupplier"
$ws.Range("M2").Value = "Parola"

# --- Add a new row 3 with the same layout/formatting as row 2 ---
$ws.Range("A2:M2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "Rakesh11"
$ws.Range("B3").Value = "us-9877"
$ws.Range("C3").Value = 7865
$ws.Range("D3").Value = "AutomationModel"
$ws.Range("E3").Value = "Ready to Deploy"
$ws.Range("F3").Value = "AutomationAsset"
$ws.Range("G3").Value = 44946
$ws.Range("H3").Value = "AutomationSupplier"
$ws.Range("I3").Value = 1234
$ws.Range("J3").Value = 35000
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = "….."
$ws.Range("M3").Value = "Parola"

# --- Sheet view: scroll so column E is leftmost, select L4 ---
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L4").Select()
